$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-update draw results: append the 2025-09-29 Pick 3 result as a new row.
$row = 13

# Column A ("2025-09-29") and column C ("250929") look numeric/date-like, so a
# plain .Value assignment would get auto-coerced into a date serial / number.
# Force text interpretation (like entering the value into a cell pre-formatted
# as Text), then clear the formatting again so the cell ends up as plain text
# with no leftover numeric format, matching the rest of the column.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-29"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "Pick 3"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "250929"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "2-7-0"

$ws.Cells.Item($row, 5).Value = "2025-09-29T21:36:03.813+04:00"
